$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.809.93'
$ws.Range("E2").Value = '  +1.12%  '

# Row 3
$ws.Range("D3").Value = '3.757.05'
$ws.Range("E3").Value = '  +0.39%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").Value = '602.77'
$ws.Range("E5").Value = '  +4.51%  '

# Row 6
$ws.Range("D6").Value = '186.89'
$ws.Range("E6").Value = '  +16.12%  '

# Row 7
$ws.Range("D7").Value = '3.751.35'
$ws.Range("E7").Value = '  +0.59%  '

# Row 8
$ws.Range("D8").Value = '0.637'
$ws.Range("E8").Value = '  -1.33%  '

# Row 9
$ws.Range("D9").Value = '0.997'
$ws.Range("E9").Value = '  -0.58%  '

# Row 10
$ws.Range("D10").Value = '0.729'
$ws.Range("E10").Value = '  +1.13%  '

# Row 11
$ws.Range("E11").Value = '  -2.28%  '

# Row 12
$ws.Range("D12").Value = '57.56'
$ws.Range("E12").Value = '  +12.57%  '

# Row 13
$ws.Range("D13").Value = '0.0000297'
$ws.Range("E13").Value = '  -3.80%  '

# Row 14
$ws.Range("D14").Value = '10.95'
$ws.Range("E14").Value = '  +1.25%  '

# Row 15
$ws.Range("D15").Value = '4.368.51'
$ws.Range("E15").Value = '  +0.46%  '

# Row 16
$ws.Range("D16").Value = '3.773.12'
$ws.Range("E16").Value = '  -0.46%  '

# Row 17
$ws.Range("D17").Value = '19.67'
$ws.Range("E17").Value = '  -2.41%  '

# Row 18
$ws.Range("D18").Value = '13.03'
$ws.Range("E18").Value = '  -2.19%  '

# Row 19
$ws.Range("D19").Value = '0.126'
$ws.Range("E19").Value = '  -1.33%  '

# Row 20
$ws.Range("E20").Value = '  -2.82%  '

# Row 21
$ws.Range("D21").Value = '69.550.77'
$ws.Range("E21").Value = '  +0.78%  '

# Row 22
$ws.Range("D22").Value = '416.39'
$ws.Range("E22").Value = '  -1.83%  '

# Row 23
$ws.Range("D23").Value = '4.66'
$ws.Range("E23").Value = '  +2.05%  '

# Row 24
$ws.Range("D24").Value = '90.02'
$ws.Range("E24").Value = '  -0.91%  '

# Row 25
$ws.Range("D25").Value = '3.08'
$ws.Range("E25").Value = '  -2.73%  '

# Row 26
$ws.Range("D26").Value = '13.03'
$ws.Range("E26").Value = '  -2.94%  '

# Row 27
$ws.Range("D27").Value = '11.20'
$ws.Range("E27").Value = '  +3.10%  '

# Row 28
$ws.Range("D28").Value = '4.00'
$ws.Range("E28").Value = '  +4.50%  '

# Row 29
$ws.Range("E29").Value = '  +2.85%  '

# Row 30
$ws.Range("D30").Value = '9.62'
$ws.Range("E30").Value = '  -4.22%  '

# Row 31
$ws.Range("D31").Value = '33.31'
$ws.Range("E31").Value = '  -1.66%  '

# Row 32
$ws.Range("D32").Value = '7.44'
$ws.Range("E32").Value = '  -4.70%  '

# Row 33
$ws.Range("D33").Value = '12.61'
$ws.Range("E33").Value = '  -3.37%  '

# Row 34
$ws.Range("E34").Value = '  -0.94%  '

# Row 35
$ws.Range("D35").Value = '44.48'
$ws.Range("E35").Value = '  -4.50%  '

# Row 36
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '614.92'
$ws.Range("E36").Value = '  +0.17%  '

# Row 37
$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '65.42'
$ws.Range("E37").Value = '  -3.16%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0913'
$ws.Range("E38").Value = '  -3.29%  '

# Row 39
$ws.Range("D39").Value = '0.410'
$ws.Range("E39").Value = '  -0.22%  '

# Row 40
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.17%  '

# Row 41
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  -0.17%  '

# Row 42
$ws.Range("E42").Value = '  -1.28%  '

# Row 43
$ws.Range("D43").Value = '3.09'
$ws.Range("E43").Value = '  -1.49%  '

# Row 44
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  +2.42%  '

# Row 45
$ws.Range("D45").Value = '3.03'
$ws.Range("E45").Value = '  -1.33%  '

# Row 46
$ws.Range("D46").Value = '0.0448'
$ws.Range("E46").Value = '  -1.32%  '

# Row 47
$ws.Range("D47").Value = '9.40'
$ws.Range("E47").Value = '  -1.85%  '

# Row 48
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -0.68%  '

# Row 49
$ws.Range("E49").Value = '  -2.31%  '

# Row 50
$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '3.23'
$ws.Range("E50").Value = '  +0.89%  '

# Row 51
$ws.Range("D51").Value = '2.800.67'
$ws.Range("E51").Value = '  +0.64%  '
